$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04558966666666667
$ws.Range("H2").Value = 0.136769
$ws.Range("I2").Value = 0.02375599288687187
$ws.Range("J2").Value = 0.02375599288687187
$ws.Range("M2").Value = 97.57717366666668
$ws.Range("N2").Value = 292.731521
$ws.Range("O2").Value = 0.3532166605548384
$ws.Range("P2").Value = 0.3532166605548384
$ws.Range("Q2").Value = 4.448510821738778
$ws.Range("R2").Value = 40.036597395649
$ws.Range("S2").Value = 0.008391012475665377
$ws.Range("T2").Value = 0.008391012475665378

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04558966666666667
$ws.Range("H3").Value = 0.136769
$ws.Range("I3").Value = 0.02375599288687187
$ws.Range("J3").Value = 0.02375599288687187
$ws.Range("O3").Value = 0.5533024543641269
$ws.Range("P3").Value = 0.5533024543641269
$ws.Range("Q3").Value = 6.968448068296333
$ws.Range("R3").Value = 62.716032614667
$ws.Range("S3").Value = 0.01314424917016295
$ws.Range("T3").Value = 0.01314424917016295

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04558966666666667
$ws.Range("H4").Value = 0.136769
$ws.Range("I4").Value = 0.02375599288687187
$ws.Range("J4").Value = 0.02375599288687187
$ws.Range("O4").Value = 0.09348088508103472
$ws.Range("P4").Value = 0.09348088508103473
$ws.Range("Q4").Value = 1.177324784894
$ws.Range("R4").Value = 10.595923064046
$ws.Range("S4").Value = 0.002220731241043548
$ws.Range("T4").Value = 0.002220731241043548

# Row 5
$ws.Range("I5").Value = 0.1978186777627204
$ws.Range("J5").Value = 0.1978186777627204
$ws.Range("M5").Value = 97.57717366666668
$ws.Range("N5").Value = 292.731521
$ws.Range("O5").Value = 0.3532166605548384
$ws.Range("P5").Value = 0.3532166605548384
$ws.Range("Q5").Value = 37.04322243907667
$ws.Range("R5").Value = 333.38900195169
$ws.Range("S5").Value = 0.06987285275472176
$ws.Range("T5").Value = 0.06987285275472177

# Row 6
$ws.Range("I6").Value = 0.1978186777627204
$ws.Range("J6").Value = 0.1978186777627204
$ws.Range("O6").Value = 0.5533024543641269
$ws.Range("P6").Value = 0.5533024543641269
$ws.Range("S6").Value = 0.1094535599251795
$ws.Range("T6").Value = 0.1094535599251795

# Row 7
$ws.Range("I7").Value = 0.1978186777627204
$ws.Range("J7").Value = 0.1978186777627204
$ws.Range("O7").Value = 0.09348088508103472
$ws.Range("P7").Value = 0.09348088508103473
$ws.Range("S7").Value = 0.0184922650828191
$ws.Range("T7").Value = 0.01849226508281911

# Row 8
$ws.Range("I8").Value = 0.7784253293504076
$ws.Range("J8").Value = 0.7784253293504078
$ws.Range("M8").Value = 97.57717366666668
$ws.Range("N8").Value = 292.731521
$ws.Range("O8").Value = 0.3532166605548384
$ws.Range("P8").Value = 0.3532166605548384
$ws.Range("Q8").Value = 145.7667342308604
$ws.Range("R8").Value = 1311.900608077743
$ws.Range("S8").Value = 0.2749527953244512
$ws.Range("T8").Value = 0.2749527953244513

# Row 9
$ws.Range("I9").Value = 0.7784253293504076
$ws.Range("J9").Value = 0.7784253293504078
$ws.Range("O9").Value = 0.5533024543641269
$ws.Range("P9").Value = 0.5533024543641269
$ws.Range("S9").Value = 0.4307046452687844
$ws.Range("T9").Value = 0.4307046452687845

# Row 10
$ws.Range("I10").Value = 0.7784253293504076
$ws.Range("J10").Value = 0.7784253293504078
$ws.Range("O10").Value = 0.09348088508103472
$ws.Range("P10").Value = 0.09348088508103473
$ws.Range("S10").Value = 0.07276788875717206
$ws.Range("T10").Value = 0.07276788875717208
